# Generate Report for Handoff
# Adds two new handed-off files (74637ef2-... and 8633d73e-...) to the
# Overview, zh-cn and de-de worksheets of the localization status report.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"
$includeReason = "Include"
$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# New item definitions, in the order they should be appended (row 6, row 7)
# ---------------------------------------------------------------------
$items = @(
  @{
    Md = "74637ef2-d2c4-47c1-853f-40bc2e87895c.md"
    ZhXlf = "74637ef2-d2c4-47c1-853f-40bc2e87895c.488e39162d363fe5ff3500b4d782a1ab26bbba82.zh-cn.xlf"
    DeXlf = "74637ef2-d2c4-47c1-853f-40bc2e87895c.488e39162d363fe5ff3500b4d782a1ab26bbba82.de-de.xlf"
    OverviewDate = "2016-26-12 06:26:47"
    ZhDate = "2016-03-12 06:26:45"
    DeDate = "2016-03-12 06:26:47"
    MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/488e39162d363fe5ff3500b4d782a1ab26bbba82/e2e/74637ef2-d2c4-47c1-853f-40bc2e87895c.md"
    ZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/488e39162d363fe5ff3500b4d782a1ab26bbba82/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/74637ef2-d2c4-47c1-853f-40bc2e87895c.488e39162d363fe5ff3500b4d782a1ab26bbba82.zh-cn.xlf"
    DeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/488e39162d363fe5ff3500b4d782a1ab26bbba82/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/74637ef2-d2c4-47c1-853f-40bc2e87895c.488e39162d363fe5ff3500b4d782a1ab26bbba82.de-de.xlf"
  },
  @{
    Md = "8633d73e-2db7-4443-9b54-cc101d84e6e3.md"
    ZhXlf = "8633d73e-2db7-4443-9b54-cc101d84e6e3.42d1603781f7447995d482e57e01eea1564f7df9.zh-cn.xlf"
    DeXlf = "8633d73e-2db7-4443-9b54-cc101d84e6e3.42d1603781f7447995d482e57e01eea1564f7df9.de-de.xlf"
    OverviewDate = "2016-26-12 06:26:47"
    ZhDate = "2016-03-12 06:26:45"
    DeDate = "2016-03-12 06:26:47"
    MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/42d1603781f7447995d482e57e01eea1564f7df9/e2e/8633d73e-2db7-4443-9b54-cc101d84e6e3.md"
    ZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42d1603781f7447995d482e57e01eea1564f7df9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8633d73e-2db7-4443-9b54-cc101d84e6e3.42d1603781f7447995d482e57e01eea1564f7df9.zh-cn.xlf"
    DeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42d1603781f7447995d482e57e01eea1564f7df9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8633d73e-2db7-4443-9b54-cc101d84e6e3.42d1603781f7447995d482e57e01eea1564f7df9.de-de.xlf"
  }
)

# ---------------------------------------------------------------------
# Sheet1 "Overview": columns A (File Name), B (zh-cn), C (de-de), D (Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$rowIdx = 6
foreach ($item in $items) {
  $wsOverview.Range("B$rowIdx").Value = $readyStatus
  $wsOverview.Range("C$rowIdx").Value = $readyStatus
  $wsOverview.Range("D$rowIdx").Value = $item.OverviewDate
  $wsOverview.Hyperlinks.Add($wsOverview.Range("A$rowIdx"), $item.MdUrl, "", "", $item.Md) | Out-Null
  $rowIdx = $rowIdx + 1
}

# ---------------------------------------------------------------------
# Sheet2 "zh-cn": columns A..K (Source File Name .. Error Detail)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$rowIdx = 6
foreach ($item in $items) {
  $wsZh.Range("C$rowIdx").Value = $readyStatus
  $wsZh.Range("E$rowIdx").Value = $item.ZhDate
  $wsZh.Range("H$rowIdx").Value = $zeroDate
  $wsZh.Range("I$rowIdx").Value = $includeReason
  $wsZh.Hyperlinks.Add($wsZh.Range("A$rowIdx"), $item.MdUrl, "", "", $item.Md) | Out-Null
  $wsZh.Hyperlinks.Add($wsZh.Range("B$rowIdx"), $item.MdUrl, "", "", ".md") | Out-Null
  $wsZh.Hyperlinks.Add($wsZh.Range("D$rowIdx"), $item.ZhXlfUrl, "", "", $item.ZhXlf) | Out-Null
  $rowIdx = $rowIdx + 1
}

# ---------------------------------------------------------------------
# Sheet3 "de-de": columns A..K (Source File Name .. Error Detail)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$rowIdx = 6
foreach ($item in $items) {
  $wsDe.Range("C$rowIdx").Value = $readyStatus
  $wsDe.Range("E$rowIdx").Value = $item.DeDate
  $wsDe.Range("H$rowIdx").Value = $zeroDate
  $wsDe.Range("I$rowIdx").Value = $includeReason
  $wsDe.Hyperlinks.Add($wsDe.Range("A$rowIdx"), $item.MdUrl, "", "", $item.Md) | Out-Null
  $wsDe.Hyperlinks.Add($wsDe.Range("B$rowIdx"), $item.MdUrl, "", "", ".md") | Out-Null
  $wsDe.Hyperlinks.Add($wsDe.Range("D$rowIdx"), $item.DeXlfUrl, "", "", $item.DeXlf) | Out-Null
  $rowIdx = $rowIdx + 1
}

Write-Output "Report generated for handoff: added $($items.Count) new files to Overview, zh-cn, de-de sheets."
